# CSU05 - Contatar Administrador: apply the edits described by the commit
# "Alteração dos casos de Uso."
#
#  1. Drop the stray _GoBack bookmark that was sitting around "COLLECTOR SHOP"
#     (left over from the previous save's last edit position).
#  2. Step 2 of the normal flow is renumbered to 3 and reworded.
#  3. Step 7 of the alternate flow is reworded ("sistema exibe mensagem" ->
#     "é exibida uma mensagem").
#  4. Step 8 of the exception flow is reworded ("mensagem de erro" -> "uma de
#     erro"), and a new _GoBack bookmark is left at the point where the user
#     stopped editing (right after "exibe uma").

$d = $word.ActiveDocument

# --- 1. Remove the old _GoBack bookmark around "COLLECTOR SHOP" -----------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- 2. Step 2 -> Step 3 wording -------------------------------------------
$d.Content.Find.Execute(
    "2 - Usuário digita seu email, escolhe a categoria, descreve sua indagação e confirma ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "3 - Usuário digita seu email, escolhe a categoria, subcategoria e descreve sua indagação confirmando.",
    2) | Out-Null

# --- 3. Step 7 wording -------------------------------------------------------
$d.Content.Find.Execute(
    "7- Usuário deixa campos sem preencher e confirma, então sistema exibe mensagem que os campos são obrigatórios e não envia mensagem;",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "7- Usuário deixa campos sem preencher e confirma, então é exibida uma mensagem que os campos são obrigatórios e não envia mensagem;",
    2) | Out-Null

# --- 4. Step 8 wording + new _GoBack bookmark -------------------------------
$d.Content.Find.Execute(
    "8- Sistema não consegue enviar mensagem e exibe mensagem de erro. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "8- Sistema não consegue enviar mensagem e exibe uma de erro. ",
    2) | Out-Null

$r = $d.Content
$r.Find.Execute(
    "8- Sistema não consegue enviar mensagem e exibe uma",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$insertPoint = $d.Range($r.End, $r.End)
$d.Bookmarks.Add("_GoBack", $insertPoint)
